# CIERRE 7 ENE 2022
# Applies the December-2021 credit postings + the stray November fix
# that were recorded in this commit.

$wb = $excel.ActiveWorkbook

$wsNov = $wb.Worksheets.Item(4)   # REMISIONES NOVIEMBRE   2021
$wsDic = $wb.Worksheets.Item(5)   # REMISIONES DICIEMBRE 2021

# ---------------------------------------------------------------------
# REMISIONES NOVIEMBRE 2021 - row 27 got its payment info completed
# ---------------------------------------------------------------------
$wsNov.Range("F27").Value = 44539
$wsNov.Range("G27").Value = 8806

# ---------------------------------------------------------------------
# REMISIONES DICIEMBRE 2021 - rows 4-13 : new credit entries for the
# December close
# ---------------------------------------------------------------------
$wsDic.Range("A4").Value = 44537
$wsDic.Range("D4").Value = "OBRADOR"
$wsDic.Range("E4").Value = 52
$wsDic.Range("F4").Value = 44537
$wsDic.Range("G4").Value = 52

$wsDic.Range("A5").Value = 44537
$wsDic.Range("D5").Value = "OBRADOR"
$wsDic.Range("E5").Value = 92
$wsDic.Range("F5").Value = 44537
$wsDic.Range("G5").Value = 92

$wsDic.Range("A6").Value = 44539
$wsDic.Range("D6").Value = "OBRADOR"
$wsDic.Range("E6").Value = 231
$wsDic.Range("F6").Value = 44544
$wsDic.Range("G6").Value = 231

$wsDic.Range("A7").Value = 44539
$wsDic.Range("D7").Value = "MAURO"
$wsDic.Range("E7").Value = 4929
$wsDic.Range("F7").Value = 44541
$wsDic.Range("G7").Value = 4929

$wsDic.Range("A8").Value = 44541
$wsDic.Range("D8").Value = "MAURO"
$wsDic.Range("E8").Value = 11506
$wsDic.Range("F8").Value = 44542
$wsDic.Range("G8").Value = 11506

$wsDic.Range("A9").Value = 44541
$wsDic.Range("D9").Value = "OBRADOR"
$wsDic.Range("E9").Value = 65
$wsDic.Range("F9").Value = 44544
$wsDic.Range("G9").Value = 65

$wsDic.Range("A10").Value = 44542
$wsDic.Range("D10").Value = "MAURO"
$wsDic.Range("E10").Value = 10389

$wsDic.Range("A11").Value = 44544
$wsDic.Range("D11").Value = "OBRADOR"
$wsDic.Range("E11").Value = 269
$wsDic.Range("F11").Value = 44544
$wsDic.Range("G11").Value = 269

$wsDic.Range("A12").Value = 44545
$wsDic.Range("D12").Value = "OBRADOR"
$wsDic.Range("E12").Value = 914
$wsDic.Range("F12").Value = 44545
$wsDic.Range("G12").Value = 914

$wsDic.Range("A13").Value = 44545
$wsDic.Range("D13").Value = "HERRADURA DAVIR"
$wsDic.Range("E13").Value = 2030

# ---------------------------------------------------------------------
# Restore on-screen selections to match where the author left off
# ---------------------------------------------------------------------
$wsNov.Select()
$wsNov.Range("H27").Select()

$wsDic.Select()
$wsDic.Range("G12").Select()
